$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 78 (pushes existing rows 78:95 down to 80:97),
# inheriting formatting (e.g. date style on column D) from the row above.
$ws.Rows.Item(78).Resize(2).Insert()

# New weekly record (2022-05-13) split by quality "Primera" / "Segunda",
# matching the existing table's column layout.
$ws.Range("A78:A79").Value2 = 1
$ws.Range("B78:B79").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C78:C79").Value2 = "Arica y Parinacota"
$ws.Range("D78:D79").Value2 = 44694
$ws.Range("E78:E79").Value2 = 15
$ws.Range("F78:F79").Value2 = 100112036
$ws.Range("G78:G79").Value2 = "Caigua"
$ws.Range("H78:H79").Value2 = "Sin especificar"
$ws.Range("O78:O79").Value2 = "Región de Arica y Parinacota"
$ws.Range("R78:R79").Value2 = "Hortaliza"

$ws.Range("I78").Value2 = "Primera"
$ws.Range("J78").Value2 = 120
$ws.Range("K78").Value2 = 13000
$ws.Range("L78").Value2 = 14000
$ws.Range("M78").Value2 = 13500
$ws.Range("N78").Value2 = "$/caja 20 kilos"
$ws.Range("P78").Value2 = 675
$ws.Range("Q78").Value2 = 20

$ws.Range("I79").Value2 = "Segunda"
$ws.Range("J79").Value2 = 130
$ws.Range("K79").Value2 = 10000
$ws.Range("L79").Value2 = 11000
$ws.Range("M79").Value2 = 10500
$ws.Range("N79").Value2 = "$/caja 20 kilos"
$ws.Range("P79").Value2 = 525
$ws.Range("Q79").Value2 = 20
